# Apply data update for 2024-07-25 (CTA violent crime YTD tracker)
$wb = $excel.ActiveWorkbook

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('D2').Value = 54
$ws.Range('F2').Value = 50
$ws.Range('H3').Value = 61
$ws.Range('C6').Value = 274
$ws.Range('D6').Value = 248
$ws.Range('F6').Value = 315
$ws.Range('G6').Value = 288
$ws.Range('H6').Value = 246
$ws.Range('I6').Value = 311
$ws.Range('K6').Value = 297
$ws.Range('C7').Value = 368
$ws.Range('D7').Value = 388
$ws.Range('F7').Value = 451
$ws.Range('G7').Value = 418
$ws.Range('H7').Value = 380
$ws.Range('I7').Value = 504
$ws.Range('K7').Value = 525

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 10
$ws.Range('K6').Value = 14
$ws.Range('K7').Value = 26

# Sheet: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('G4').Value = 1
$ws.Range('G5').Value = 3

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K31').Value = 26
$ws.Range('D34').Value = 5
$ws.Range('F46').Value = 11
$ws.Range('G46').Value = 15
$ws.Range('D52').Value = 54
$ws.Range('F52').Value = 45
$ws.Range('D53').Value = 4
$ws.Range('K60').Value = 2
$ws.Range('K63').Value = 5
$ws.Range('C64').Value = 12
$ws.Range('H66').Value = 3
$ws.Range('K73').Value = 11
$ws.Range('K75').Value = 17
$ws.Range('H76').Value = 13
$ws.Range('I76').Value = 27
$ws.Range('H79').Value = 3
$ws.Range('G93').Value = 3
$ws.Range('C97').Value = 368
$ws.Range('D97').Value = 388
$ws.Range('F97').Value = 451
$ws.Range('G97').Value = 418
$ws.Range('H97').Value = 380
$ws.Range('I97').Value = 504
$ws.Range('K97').Value = 525

# Sheet: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('D6').Value = 31
$ws.Range('F6').Value = 35
$ws.Range('D7').Value = 54
$ws.Range('F7').Value = 45

# Sheet: North Center
$ws = $wb.Worksheets.Item('North Center')
$ws.Range('E3').Value = 3
$ws.Range('E4').Value = 5

# Sheet: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('D2').Value = 1
$ws.Range('D5').Value = 4
$ws.Range('D6').Value = 5

# Sheet: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K5').Value = 7
$ws.Range('K6').Value = 17

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('C5').Value = 10
$ws.Range('C6').Value = 12

# Sheet: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K5').Value = 8
$ws.Range('K6').Value = 11

# Sheet: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('H3').Value = 2
$ws.Range('I6').Value = 15
$ws.Range('H7').Value = 13
$ws.Range('I7').Value = 27

# Sheet: Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('D4').Value = 3
$ws.Range('D5').Value = 4

# Sheet: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('F2').Value = 2
$ws.Range('G5').Value = 8
$ws.Range('F6').Value = 11
$ws.Range('G6').Value = 15

# Sheet: West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('F5').Value = 2
$ws.Range('F6').Value = 3

# Sheet: Norwood Park
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('D5').Value = 3
$ws.Range('D6').Value = 3
